$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unhide the "discounts" block (rows 24-27) before touching cell values, so that
# Excel does not recompute/auto-fit the row height of row 24 when its value changes. ---
$ws.Rows.Item(24).Hidden = $false
$ws.Rows.Item(25).Hidden = $false
$ws.Rows.Item(26).Hidden = $false
$ws.Rows.Item(27).Hidden = $false

# --- Row 22: add the missing "getProductsFromCart" label under the cart block ---
$ws.Range("C22").Value = "getProductsFromCart"

# --- Row 24: rename "coupons" / "addCoupon" block to "discounts" / "addDiscount" ---
$ws.Range("B24").Value = "discounts"
$ws.Range("C24").Value = "addDiscount"

# --- View changes: zoom level and selected cell ---
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("L10").Select()

# --- Page setup: print scale ---
$ws.PageSetup.Zoom = 108
